$wb = $excel.ActiveWorkbook

# --- Notes sheet: insert the new note line about the wide-format sheet ---
$notes = $wb.Worksheets.Item("Notes")
$notes.Rows.Item(8).Insert()
$notes.Cells.Item(8, 1).Value = "On the 'Data-wide-value' sheet, we have provided the indicator in a wide format. The values you see listed there are from the 'value' column."
$notes.Rows.Item(11).Insert()

# --- Add the new Data-wide-value sheet after the Data sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wide = $wb.Worksheets.Add($null, $lastSheet)
$wide.Name = "Data-wide-value"

$wide.Cells.Item(1, 1).Value = "id"
$wide.Range("B1:C1").NumberFormat = "@"
$wide.Cells.Item(1, 2).Value = "2013"
$wide.Cells.Item(1, 3).Value = "2014"

$wide.Cells.Item(2, 1).Value = "ET"
$wide.Cells.Item(2, 2).Value = 3359031639
$wide.Cells.Item(2, 3).Value = 3565493249

$wide.Cells.Item(3, 1).Value = "NG"
$wide.Cells.Item(3, 2).Value = 0
$wide.Cells.Item(3, 3).Value = 2545783781

# --- Restore Notes as the active/selected sheet ---
$notes.Activate()
